$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 269   # F2: 267 -> 269
$ws1.Cells.Item(6, 6).Value = 1076   # F6: 1075 -> 1076
$ws1.Cells.Item(7, 6).Value = 0   # F7: 1406 -> 0
$ws1.Cells.Item(8, 6).Value = 579   # F8: 578 -> 579
$ws1.Cells.Item(10, 6).Value = 739   # F10: 737 -> 739
$ws1.Cells.Item(12, 6).Value = 133   # F12: 131 -> 133
$ws1.Cells.Item(14, 6).Value = 417   # F14: 415 -> 417
$ws1.Cells.Item(15, 6).Value = 1313   # F15: 1308 -> 1313
$ws1.Cells.Item(16, 6).Value = 99   # F16: 98 -> 99
$ws1.Cells.Item(18, 6).Value = 266   # F18: 265 -> 266
$ws1.Cells.Item(23, 6).Value = 10   # F23: 9 -> 10
$ws1.Cells.Item(24, 6).Value = 5630   # F24: 5618 -> 5630
$ws1.Cells.Item(25, 6).Value = 53   # F25: 52 -> 53
$ws1.Cells.Item(27, 6).Value = 89   # F27: 87 -> 89
$ws1.Cells.Item(29, 6).Value = 14168   # F29: 14148 -> 14168
$ws1.Cells.Item(30, 6).Value = 1412   # F30: 1410 -> 1412
$ws1.Cells.Item(32, 6).Value = 90   # F32: 89 -> 90
$ws1.Cells.Item(33, 6).Value = 83   # F33: 82 -> 83
$ws1.Cells.Item(34, 6).Value = 423   # F34: 419 -> 423
$ws1.Cells.Item(35, 6).Value = 582   # F35: 578 -> 582
$ws1.Cells.Item(36, 6).Value = 4169   # F36: 4166 -> 4169
$ws1.Cells.Item(37, 6).Value = 115   # F37: 107 -> 115

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 269   # F2: 267 -> 269
$ws4.Cells.Item(6, 6).Value = 1076   # F6: 1075 -> 1076
$ws4.Cells.Item(7, 6).Value = 1408   # F7: 1406 -> 1408
$ws4.Cells.Item(8, 6).Value = 0   # F8: 578 -> 0
$ws4.Cells.Item(10, 6).Value = 739   # F10: 737 -> 739
$ws4.Cells.Item(12, 6).Value = 133   # F12: 131 -> 133
$ws4.Cells.Item(14, 6).Value = 417   # F14: 415 -> 417
$ws4.Cells.Item(15, 6).Value = 1313   # F15: 1308 -> 1313
$ws4.Cells.Item(16, 6).Value = 99   # F16: 98 -> 99
$ws4.Cells.Item(18, 6).Value = 266   # F18: 265 -> 266
$ws4.Cells.Item(25, 6).Value = 10   # F25: 9 -> 10
$ws4.Cells.Item(27, 6).Value = 5630   # F27: 5618 -> 5630
$ws4.Cells.Item(28, 6).Value = 53   # F28: 52 -> 53
$ws4.Cells.Item(30, 6).Value = 89   # F30: 87 -> 89
$ws4.Cells.Item(32, 6).Value = 14168   # F32: 14148 -> 14168
$ws4.Cells.Item(33, 6).Value = 1412   # F33: 1410 -> 1412
$ws4.Cells.Item(35, 6).Value = 90   # F35: 89 -> 90
$ws4.Cells.Item(36, 6).Value = 83   # F36: 82 -> 83
$ws4.Cells.Item(37, 6).Value = 423   # F37: 419 -> 423
$ws4.Cells.Item(38, 6).Value = 582   # F38: 578 -> 582
$ws4.Cells.Item(39, 6).Value = 4169   # F39: 4166 -> 4169
$ws4.Cells.Item(40, 6).Value = 115   # F40: 107 -> 115
